# 2.03 Changed mosfet to the much more robust BTS133
# Replace the IRLB3034PbF HEXFET MOSFET BOM line with the BTS133BKSA1
# smart low-side power switch, and log the change on the Revision sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # electronicloadBOM
$ws2 = $wb.Worksheets.Item(2)   # Revision

# ---------------------------------------------------------------------
# 1. BOM sheet, row 15 (the Q1 mosfet line)
# ---------------------------------------------------------------------

# Value (B) - wraps, new Calculation style
$ws1.Range("B15").Value = "BTS133BKSA1"
$ws1.Range("B15").Style = "Calculation"
$ws1.Range("B15").WrapText = $true

# Package / Parts (C/D) stay the same - untouched on purpose
$ws1.Range("C15").Value = "TO220BV"
$ws1.Range("D15").Value = "Q1"

# Description (E) - new Calculation style, no special alignment
$ws1.Range("E15").Value = "Smart low side power switch"
$ws1.Range("E15").Style = "Calculation"

# Tolerance (F) stays blank, but now styled
$ws1.Range("F15").Style = "Calculation"

# Manufacturer (G) / MPN (H) - new Calculation style, left aligned
$ws1.Range("G15").Value = "Infineon"
$ws1.Range("G15").Style = "Calculation"
$ws1.Range("G15").HorizontalAlignment = -4131

$ws1.Range("H15").Value = "BTS133BKSA1"
$ws1.Range("H15").Style = "Calculation"
$ws1.Range("H15").HorizontalAlignment = -4131

# Ebay (I) - no more a genuine listing, replace the hyperlink with a
# plain warning note instead.
$ws1.Range("I15").Value = "There's are fake ones on ebay. Please buy from Digikey, Mouser, Farnell etc. "
$ws1.Range("I15").Style = "Calculation"

# Remove the now-stale I15 hyperlink (to the fake IRLB3034 listing) while
# keeping every other hyperlink on the sheet intact. The host only
# supports wiping the whole collection in one shot, so rebuild everything
# except I15 (same cells/URLs/order as before -> same rIds as before,
# just shifted down by one).
$links = @(
    @("I13", "http://www.ebay.com/itm/10PCS-STM32F103C8T6-ARM-STM32-Minimum-System-Development-Board-Module-Arduino/171907093814"),
    @("I6",  "http://www.ebay.com/itm/100pcs-0603-SMD-Resistor-1-8K-ohm-1K8-Tol-5-RoHS-1-10W/291817476969"),
    @("I20", "http://www.ebay.com/itm/3-Sheets-CE-Certificated-Label-Stickers-Adhesive-Sticker-Markers-Home-Appliance/172533333629"),
    @("I17", "http://www.ebay.com/itm/20pcs-2-54mm-Pitch-1X-20-Pin-Female-Single-Row-Straight-Header-Connector-PCB-DIY/132035408776"),
    @("I10", "http://www.ebay.com/itm/200Pcs-0603-SMD-Resistor-Resistors-1K-910K-Ohm-Ω-1-High-Quality-Free-Shipping/192101507976"),
    @("I11", "http://www.ebay.com/itm/200Pcs-0603-SMD-Resistor-Resistors-1K-910K-Ohm-Ω-1-High-Quality-Free-Shipping/192101507976"),
    @("I8",  "http://www.ebay.com/itm/200Pcs-0603-SMD-Resistor-Resistors-1K-910K-Ohm-Ω-1-High-Quality-Free-Shipping/192101507976"),
    @("I3",  "http://www.ebay.com/itm/200Pcs-0603-SMD-Resistor-Resistors-1K-910K-Ohm-Ω-1-High-Quality-Free-Shipping/192101507976"),
    @("I5",  "http://www.ebay.com/itm/200PCS-150R-150-ohm-Ω-1-1-10W-SMD-Chip-Resistor-0603-1-6mm-0-8mm/231949532343"),
    @("I9",  "http://www.ebay.com/itm/100PCS-22nF-223-10-50V-X7R-0603-1608-SMD-capacitor-MLCC/232244234398"),
    @("I4",  "http://www.ebay.com/itm/100PCS-0603-1608-100n-100nF-104K-10-50V-X7R-SMD-capacitor-MLCC/391284797144"),
    @("I7",  "http://www.ebay.com/itm/100PCS-1uF-105K-10-16V-X7R-0603-1608-SMD-capacitor-MLCC-1-6mm-0-8mm/232352831873"),
    @("I18", "http://www.ebay.com/itm/100PCS-Self-Adhesive-Rubber-Feet-Clear-Semicircle-Bumpers-Door-Buffer-Pad/361356150996"),
    @("I2",  "http://www.ebay.com/itm/100PCS-2512-6432-1W-0-3-ohm-R300-0-3R-1-SMD-resistors/291647639892"),
    @("I14", "http://www.ebay.com/itm/20pcs-2Pin-Plug-in-Terminal-Block-DG128-Screw-KF128-2P-Pitch-5-08MM-300V-10A/391527936166"),
    @("I19", "http://www.ebay.com/itm/20-sets-TO-220-Heatsink-Mounting-Insulator-Kits/190888986454"),
    @("I16", "http://www.ebay.com/itm/20pc-Aluminum-Heatsink-Ak-122-H-45mm-Size-42x25x45mm-Color-Natural-TO-3P-TO-220/131284780674")
)
$ws1.Hyperlinks.Delete()
foreach ($link in $links) {
    $ws1.Hyperlinks.Add($ws1.Range($link[0]), $link[1])
}

# ---------------------------------------------------------------------
# 2. Revision sheet - log the change as 2.03
# ---------------------------------------------------------------------
$ws2.Range("B6").Value = 2.03
$ws2.Range("C6").Value = "Replaced mosfet with BTS133BKSA1"
$ws2.Range("B6:C6").Style = "Calculation"
$ws2.Range("B6:C6").HorizontalAlignment = -4131

# ---------------------------------------------------------------------
# 3. Restore selections (Revision first, BOM last, so the BOM sheet
#    stays the active tab like in the original file).
# ---------------------------------------------------------------------
$ws2.Range("C9").Select()
$ws1.Range("D23").Select()
